$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 88
$ws.Range("F2").Value = 0.8627450980392157
$ws.Range("G2").Value = 0.8627450980392157
$ws.Range("H2").Value = 0.09737246067906356
$ws.Range("I2").Value = 0.08400761313487837
$ws.Range("J2").Value = 462691.8401053585
$ws.Range("K2").Value = 168731.9202536792
$ws.Range("M2").Value = 168731.9202536792
$ws.Range("N2").Value = 631423.7603590377
$ws.Range("O2").Value = 10131360.5088
$ws.Range("P2").Value = 9723619.5787
$ws.Range("Q2").Value = 0.01665441873350774
$ws.Range("R2").Value = 0.01735278914276877

# Row 3
$ws.Range("E3").Value = 88
$ws.Range("F3").Value = 0.8543689320388349
$ws.Range("G3").Value = 0.8543689320388349
$ws.Range("H3").Value = 0.09684574533725648
$ws.Range("I3").Value = 0.08274199601629682
$ws.Range("J3").Value = 477113.6797688863
$ws.Range("K3").Value = 173344.4200914731
$ws.Range("M3").Value = 173344.4200914731
$ws.Range("N3").Value = 650458.0998603592
$ws.Range("O3").Value = 10494911.028964
$ws.Range("P3").Value = 10087537.870961
$ws.Range("Q3").Value = 0.0165169975822639
$ws.Range("R3").Value = 0.01718401678475774

# Row 4
$ws.Range("E4").Value = 88
$ws.Range("F4").Value = 0.8461538461538461
$ws.Range("G4").Value = 0.8461538461538461
$ws.Range("H4").Value = 0.0971597766805547
$ws.Range("I4").Value = 0.08221211872970012
$ws.Range("J4").Value = 504436.7308120827
$ws.Range("K4").Value = 179919.5730192822
$ws.Range("M4").Value = 179919.5730192822
$ws.Range("N4").Value = 684356.3038313651
$ws.Range("O4").Value = 10907997.05983292
$ws.Range("P4").Value = 10499652.70708983
$ws.Range("Q4").Value = 0.01649428140036903
$ws.Range("R4").Value = 0.01713576420463818

# Row 5
$ws.Range("E5").Value = 89
$ws.Range("F5").Value = 0.8476190476190476
$ws.Range("G5").Value = 0.8476190476190476
$ws.Range("H5").Value = 0.0964615420761171
$ws.Range("I5").Value = 0.08176264042642307
$ws.Range("J5").Value = 524764.3891885336
$ws.Range("K5").Value = 187509.438435905
$ws.Range("M5").Value = 187509.438435905
$ws.Range("N5").Value = 712273.8276244387
$ws.Range("O5").Value = 11375845.10502791
$ws.Range("P5").Value = 10965150.42170252
$ws.Range("Q5").Value = 0.01648312162346773
$ws.Range("R5").Value = 0.01710048938907224

# Row 6
$ws.Range("E6").Value = 91
$ws.Range("F6").Value = 0.8584905660377359
$ws.Range("G6").Value = 0.8584905660377359
$ws.Range("H6").Value = 0.0955045359716872
$ws.Range("I6").Value = 0.08198974314550507
$ws.Range("J6").Value = 548908.0304184185
$ws.Range("K6").Value = 196520.0763660966
$ws.Range("M6").Value = 196520.0763660966
$ws.Range("N6").Value = 745428.1067845151
$ws.Range("O6").Value = 11767757.80777875
$ws.Range("P6").Value = 11353292.2839536
$ws.Range("Q6").Value = 0.01669987431558053
$ws.Range("R6").Value = 0.0173095232159091

# Row 7
$ws.Range("D7").Value = 101
$ws.Range("E7").Value = 87
$ws.Range("F7").Value = 0.8613861386138614
$ws.Range("G7").Value = 0.8529411764705882
$ws.Range("H7").Value = 0.09777407078138732
$ws.Range("I7").Value = 0.08339553096059507
$ws.Range("J7").Value = 459306.928978237
$ws.Range("K7").Value = 167039.4646901185
$ws.Range("M7").Value = 167039.4646901185
$ws.Range("N7").Value = 626346.3936683555
$ws.Range("O7").Value = 10084600.3388
$ws.Range("P7").Value = 9676859.4087
$ws.Range("Q7").Value = 0.01656381602426448
$ws.Range("R7").Value = 0.01726174346812782

# Row 8
$ws.Range("D8").Value = 102
$ws.Range("E8").Value = 88
$ws.Range("F8").Value = 0.8627450980392157
$ws.Range("G8").Value = 0.8543689320388349
$ws.Range("H8").Value = 0.09815169676873819
$ws.Range("I8").Value = 0.08385776034610644
$ws.Range("J8").Value = 484236.7288196762
$ws.Range("K8").Value = 176905.9446168681
$ws.Range("M8").Value = 176905.9446168681
$ws.Range("N8").Value = 661142.6734365443
$ws.Range("O8").Value = 10481710.504064
$ws.Range("P8").Value = 10074337.346061
$ws.Range("Q8").Value = 0.01687758353450781
$ws.Range("R8").Value = 0.01756005765342344

# Row 9
$ws.Range("E9").Value = 88
$ws.Range("F9").Value = 0.8461538461538461
$ws.Range("G9").Value = 0.8461538461538461
$ws.Range("H9").Value = 0.09881066311551624
$ws.Range("I9").Value = 0.08360902263620605
$ws.Range("J9").Value = 510062.6018105842
$ws.Range("K9").Value = 182732.508518533
$ws.Range("M9").Value = 182732.508518533
$ws.Range("N9").Value = 692795.1103291172
$ws.Range("O9").Value = 10837517.91078592
$ws.Range("P9").Value = 10429173.55804283
$ws.Range("Q9").Value = 0.01686110325471023
$ws.Range("R9").Value = 0.0175212836857636

# Row 10
$ws.Range("E10").Value = 89
$ws.Range("F10").Value = 0.8476190476190476
$ws.Range("G10").Value = 0.8476190476190476
$ws.Range("H10").Value = 0.09810748281467674
$ws.Range("I10").Value = 0.08315777114767837
$ws.Range("J10").Value = 528942.5728075609
$ws.Range("K10").Value = 189598.5302454186
$ws.Range("M10").Value = 189598.5302454186
$ws.Range("N10").Value = 718541.1030529796
$ws.Range("O10").Value = 11253820.6075095
$ws.Range("P10").Value = 10843125.92418412
$ws.Range("Q10").Value = 0.01684748112289106
$ws.Range("R10").Value = 0.01748559701059497

# Row 11
$ws.Range("E11").Value = 89
$ws.Range("F11").Value = 0.839622641509434
$ws.Range("G11").Value = 0.839622641509434
$ws.Range("H11").Value = 0.09785518137528718
$ws.Range("I11").Value = 0.08216142587170339
$ws.Range("J11").Value = 546463.8167449427
$ws.Range("K11").Value = 195297.9695293586
$ws.Range("M11").Value = 195297.9695293586
$ws.Range("N11").Value = 741761.7862743011
$ws.Range("O11").Value = 11744042.02153478
$ws.Range("P11").Value = 11329576.49770964
$ws.Range("Q11").Value = 0.01662953599546435
$ws.Range("R11").Value = 0.01723788789182363

# Row 12
$ws.Range("E12").Value = 87
$ws.Range("F12").Value = 0.8529411764705882
$ws.Range("G12").Value = 0.8529411764705882
$ws.Range("H12").Value = 0.105155020195745
$ws.Range("I12").Value = 0.08969104663754716
$ws.Range("J12").Value = 519414.7856197282
$ws.Range("K12").Value = 197093.3930108641
$ws.Range("M12").Value = 197093.3930108641
$ws.Range("N12").Value = 716508.1786305922
$ws.Range("O12").Value = 10094300.9488
$ws.Range("P12").Value = 9686560.0187
$ws.Range("Q12").Value = 0.01952521467415674
$ws.Range("R12").Value = 0.02034709872548907

# Row 13
$ws.Range("D13").Value = 102
$ws.Range("E13").Value = 85
$ws.Range("F13").Value = 0.8333333333333334
$ws.Range("G13").Value = 0.8252427184466019
$ws.Range("H13").Value = 0.1124218121923839
$ws.Range("I13").Value = 0.09277528190633627
$ws.Range("J13").Value = 589272.077115087
$ws.Range("K13").Value = 229423.6187645734
$ws.Range("M13").Value = 229423.6187645734
$ws.Range("N13").Value = 818695.6958796604
$ws.Range("O13").Value = 10505163.513564
$ws.Range("P13").Value = 10097790.355561
$ws.Range("Q13").Value = 0.02183912877399267
$ws.Range("R13").Value = 0.02272018042424762

# Row 14
$ws.Range("E14").Value = 84
$ws.Range("F14").Value = 0.8076923076923077
$ws.Range("G14").Value = 0.8076923076923077
$ws.Range("H14").Value = 0.1180246687376596
$ws.Range("I14").Value = 0.09532761705734048
$ws.Range("J14").Value = 649734.5660466086
$ws.Range("K14").Value = 252568.4906365452
$ws.Range("M14").Value = 252568.4906365452
$ws.Range("N14").Value = 902303.0566831537
$ws.Range("O14").Value = 10901225.88937092
$ws.Range("P14").Value = 10492881.53662783
$ws.Range("Q14").Value = 0.02316881543412548
$ws.Range("R14").Value = 0.02407046050743035

# Row 15
$ws.Range("H15").Value = 0.1126370358327302
$ws.Range("I15").Value = 0.09654603071376873
$ws.Range("J15").Value = 692419.1877078008
$ws.Range("K15").Value = 271336.8376955385
$ws.Range("M15").Value = 271336.8376955385
$ws.Range("N15").Value = 963756.0254033392
$ws.Range("O15").Value = 11465708.21445205
$ws.Range("P15").Value = 11055013.53112666
$ws.Range("Q15").Value = 0.0236650743783563
$ws.Range("R15").Value = 0.02454423388370882

# Row 16
$ws.Range("E16").Value = 90
$ws.Range("F16").Value = 0.8490566037735849
$ws.Range("G16").Value = 0.8490566037735849
$ws.Range("H16").Value = 0.1130249947913053
$ws.Range("I16").Value = 0.09596461821903285
$ws.Range("J16").Value = 711267.298988305
$ws.Range("K16").Value = 277699.7106510397
$ws.Range("M16").Value = 277699.7106510397
$ws.Range("N16").Value = 988967.0096393448
$ws.Range("O16").Value = 11792951.18548561
$ws.Range("P16").Value = 11378485.66166047
$ws.Range("Q16").Value = 0.02354794031478937
$ws.Range("R16").Value = 0.0244056826987744
